$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 34, shifting existing rows 34:106 down to 35:107
$ws.Rows.Item(34).Insert()

# Populate the newly inserted row 34 with the new record's data
$ws.Cells.Item(34, 1).Value = 10
$ws.Cells.Item(34, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(34, 3).Value = "La Araucanía"
$ws.Cells.Item(34, 4).Value = 44526
$ws.Cells.Item(34, 5).Value = 9
$ws.Cells.Item(34, 6).Value = 100112012
$ws.Cells.Item(34, 7).Value = "Espinaca"
$ws.Cells.Item(34, 8).Value = "Sin especificar"
$ws.Cells.Item(34, 9).Value = "Primera"
$ws.Cells.Item(34, 10).Value = 40
$ws.Cells.Item(34, 11).Value = 8000
$ws.Cells.Item(34, 12).Value = 8000
$ws.Cells.Item(34, 13).Value = 8000
$ws.Cells.Item(34, 14).Value = "$/docena de atados"
$ws.Cells.Item(34, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(34, 16).Value = 2667
$ws.Cells.Item(34, 17).Value = 3
$ws.Cells.Item(34, 18).Value = "Hortaliza"

# Ensure the date cell keeps the same date-number format used by the rest of column D
$ws.Cells.Item(34, 4).NumberFormat = $ws.Cells.Item(35, 4).NumberFormat
